$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert 3 new rows starting at row 60. This pushes the existing
# "GET /myServices/:token" block (rows 60-62) down to rows 63-65,
# and everything below shifts down by 3 as well, preserving all
# existing formatting/styles automatically.
$ws.Rows("60:62").Insert(-4121)

# Populate the newly inserted rows with the new "GET /id=:id" (Service by ID) endpoint.
$ws.Range("B60").Value = "GET"
$ws.Range("C60").Value = "/id=:id"
$ws.Range("E60").Value = 404
$ws.Range("F60").Value = "service not found"

$ws.Range("E61").Value = 200
$ws.Range("G61").Value = "Service"

# Update the view so it matches the edited area.
$ws.Range("F62").Select()
